$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: OTHER -> UNKNOWN
$ws.Range("D1").Value = "UNKNOWN"

# Row 2
$ws.Range("A2").Value = $False
$ws.Range("B2").Value = $True
$ws.Range("C2").Value = $False
$ws.Range("D2").Value = $False
$ws.Range("E2").Value = 2273

# Row 3
$ws.Range("A3").Value = $False
$ws.Range("B3").Value = $False
$ws.Range("C3").Value = $True
$ws.Range("D3").Value = $False
$ws.Range("E3").Value = 1994

# Row 4
$ws.Range("A4").Value = $True
$ws.Range("B4").Value = $False
$ws.Range("C4").Value = $False
$ws.Range("D4").Value = $False
$ws.Range("E4").Value = 906

# Row 5
$ws.Range("A5").Value = $False
$ws.Range("B5").Value = $True
$ws.Range("C5").Value = $True
$ws.Range("D5").Value = $False
$ws.Range("E5").Value = 575

# Row 6
$ws.Range("A6").Value = $True
$ws.Range("B6").Value = $False
$ws.Range("C6").Value = $True
$ws.Range("D6").Value = $False
$ws.Range("E6").Value = 518

# Row 7
$ws.Range("A7").Value = $True
$ws.Range("B7").Value = $True
$ws.Range("C7").Value = $True
$ws.Range("D7").Value = $False
$ws.Range("E7").Value = 300

# Row 8
$ws.Range("A8").Value = $True
$ws.Range("B8").Value = $True
$ws.Range("C8").Value = $False
$ws.Range("D8").Value = $False
$ws.Range("E8").Value = 190

# Row 9
$ws.Range("A9").Value = $False
$ws.Range("B9").Value = $False
$ws.Range("C9").Value = $False
$ws.Range("D9").Value = $True
$ws.Range("E9").Value = 149

# Row 10
$ws.Range("A10").Value = $False
$ws.Range("B10").Value = $True
$ws.Range("C10").Value = $False
$ws.Range("D10").Value = $True
$ws.Range("E10").Value = 105

# Row 11
$ws.Range("A11").Value = $False
$ws.Range("B11").Value = $False
$ws.Range("C11").Value = $True
$ws.Range("D11").Value = $True
$ws.Range("E11").Value = 40

# Row 12
$ws.Range("A12").Value = $False
$ws.Range("B12").Value = $True
$ws.Range("C12").Value = $True
$ws.Range("D12").Value = $True
$ws.Range("E12").Value = 33

# Row 13
$ws.Range("A13").Value = $True
$ws.Range("B13").Value = $False
$ws.Range("C13").Value = $True
$ws.Range("D13").Value = $True
$ws.Range("E13").Value = 24

# Row 14
$ws.Range("A14").Value = $True
$ws.Range("B14").Value = $True
$ws.Range("C14").Value = $True
$ws.Range("D14").Value = $True
$ws.Range("E14").Value = 23

# Row 15
$ws.Range("A15").Value = $True
$ws.Range("B15").Value = $True
$ws.Range("C15").Value = $False
$ws.Range("D15").Value = $True
$ws.Range("E15").Value = 13

# Row 16
$ws.Range("A16").Value = $True
$ws.Range("B16").Value = $False
$ws.Range("C16").Value = $False
$ws.Range("D16").Value = $True
$ws.Range("E16").Value = 9
